$d = $word.ActiveDocument

# 1) Simple bullet-prefix / figure-fix replacements via Find & Replace
$replacements = @(
    @("Order intake increased by 53% to 514 MSEK (335).", "• Order intake increased by 53% to 514 MSEK (335)."),
    @("Net sales rose by 390% to 357 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cbb53b0a1b2fadf4e040cdc2155a7340de24aca93cba53b0a1b2fadf4e040cdc2155a7340de24aca93cbb53b0a1b2fadf4e040cdc2155a7340de24aca93cbfadf53b0a1b2fadf4e040cdc2155a7340de24aca93cbe53b0a1b2fadf4e040cdc2155a7340de24aca93cbcdc53b0a1b2fadf4e040cdc2155a7340de24aca93cba53b0a1b2fadf4e040cdc2155a7340de24aca93cbde53b0a1b2fadf4e040cdc2155a7340de24aca93cbaca53b0a1b2fadf4e040cdc2155a7340de24aca93cbcb).", "• Net sales rose by 390% to 357 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb)."),
    @("EBIT amounted to 9,8,7,9 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb). ", "• EBIT amounted to 9,8 MSEK (7,9). "),
    @("Profit after tax of 1,9,-0,8 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb).", "• Profit after tax of 1,9 MSEK (-0,8)."),
    @("Earnings per share were 0,07,-0,03 SEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb).", "• Earnings per share were 0,07 SEK (-0,03)."),
    @("Order intake increased by 26% to 932 MSEK (741).", "• Order intake increased by 26% to 932 MSEK (741)."),
    @("Net sales of 741 MSEK (753).", "• Net sales of 741 MSEK (753)."),
    @("EBIT amounted to 14,5,16,3 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb). ", "• EBIT amounted to 14,5 MSEK (16,3). "),
    @("Profit after tax of 3,6,4,4 MSEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb).", "• Profit after tax of 3,6 MSEK (4,4)."),
    @("Earnings per share were 0,12,0,16 SEK (53b0a1b2fadf4e040cdc2155a7340de24aca93cb,53b0a1b2fadf4e040cdc2155a7340de24aca93cb).", "• Earnings per share were 0,12 SEK (0,16)."),
    @("Acquisition of ELTEC of Germany.", "• Acquisition of ELTEC of Germany."),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 2) Merge the "We think…" paragraph, the blank separator after it, and the
#    "EBIT was up…" paragraph into a single paragraph with new, edited text
#    (the blank paragraph that used to follow "EBIT was up…" stays put).
$mergedText = "We think that the loss of sales resulting from the component shortage that have restricted shipments was 2453b0a1b2fadf4e040cdc2155a7340de24aca93cb% of sales in the second quarter, or nearly 53b0a1b2fadf4e040cdc2155a7340de24aca93cb MSEK. The Group was still able to increase sales by 53b0a1b2fadf4e040cdc2155a7340de24aca93cb% to almost 53b0a1b2fadf4e040cdc2155a7340de24aca93cb MSEK in the period. In this context, I’d like to acknowledge how my co-workers have succeeded in dealing with the challenges of the component shortage through flexibility, finding new solutions and re-planning.EBIT was up by 53b0a1b2fadf4e040cdc2155a7340de24aca93cb%, but negatively impacted by the loss of sales. As CEO, I can state that earnings are too low, and that there is clearly some way to go to achieve our profitability targets. "

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("We think that the loss of sales")) {
        $target = $i
        break
    }
}

$p = $d.Paragraphs($target)
$p.Range.Text = $mergedText

# Remove the next two paragraphs (blank separator, old EBIT paragraph) that used to
# follow; the blank paragraph that came after the EBIT paragraph is left intact.
for ($k = 0; $k -lt 2; $k++) {
    $d.Paragraphs($target + 1).Range.Delete() | Out-Null
}

